$wb = $excel.ActiveWorkbook

# 1. "Not yet handed off" -> "Handed back" everywhere (Overview, zh-cn, de-de all
#    share the same underlying text, so a global sweep keeps every sheet synced
#    the same way the shared-string table does in the real workbook).
foreach ($sheet in $wb.Worksheets) {
    $used = $sheet.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        for ($c = 1; $c -le $used.Columns.Count; $c++) {
            $cell = $sheet.Cells.Item($r, $c)
            if ($cell.Value() -eq "Not yet handed off") {
                $cell.Value = "Handed back"
            }
        }
    }
}

# 2. For the zh-cn / de-de report sheets: the handback step fills in the
#    "Latest Target File" (col E) and "Latest Handback File" (col F) columns
#    (same links as "Source File Name" / "Latest Handoff File") for the two
#    data rows, and stamps "Latest Handback DateTime" (col G).
function Update-ReportSheet($sheetName, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    $addrs = @{}
    $displays = @{}
    foreach ($h in $ws.Hyperlinks) {
        $a1 = $h.Range.Address()
        $addrs[$a1] = $h.Address
        $displays[$a1] = $h.TextToDisplay
    }

    $ws.Hyperlinks.Add($ws.Range("E2"), $addrs['$A$2'], $null, $null, $displays['$A$2'])
    $ws.Hyperlinks.Add($ws.Range("F2"), $addrs['$C$2'], $null, $null, $displays['$C$2'])
    $ws.Hyperlinks.Add($ws.Range("E3"), $addrs['$A$3'], $null, $null, $displays['$A$3'])
    $ws.Hyperlinks.Add($ws.Range("F3"), $addrs['$C$3'], $null, $null, $displays['$C$3'])

    # Match the workbook's existing "HyperLink" look (underline + the custom
    # blue FF6495ED) instead of Excel's theme-based default hyperlink style.
    foreach ($addr in @("E2", "F2", "E3", "F3")) {
        $rng = $ws.Range($addr)
        $rng.Font.Underline = 2
        $rng.Font.Color = 15570276
        $rng.Font.Name = "Calibri"
        $rng.Font.Size = 11
    }

    $ws.Range("G2").Value = $handbackDateTime
    $ws.Range("G3").Value = $handbackDateTime
}

Update-ReportSheet "zh-cn" "2016-01-07 04:07:03"
Update-ReportSheet "de-de" "2016-01-07 04:07:27"
